$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 57.14035266666667
$ws.Range("H2").Value = 171.421058
$ws.Range("I2").Value = 0.7274038390747541
$ws.Range("J2").Value = 0.7274038390747541
$ws.Range("M2").Value = 4.621579
$ws.Range("N2").Value = 13.864737
$ws.Range("O2").Value = 0.1778708528171788
$ws.Range("P2").Value = 0.1778708528171788
$ws.Range("Q2").Value = 264.0786539368607
$ws.Range("R2").Value = 2376.707885431746
$ws.Range("S2").Value = 0.1293839411987164
$ws.Range("T2").Value = 0.1293839411987164

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 57.14035266666667
$ws.Range("H3").Value = 171.421058
$ws.Range("I3").Value = 0.7274038390747541
$ws.Range("J3").Value = 0.7274038390747541
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5971062807549863
$ws.Range("P3").Value = 0.5971062807549863
$ws.Range("Q3").Value = 886.5028776867309
$ws.Range("R3").Value = 7978.525899180579
$ws.Range("S3").Value = 0.434337400956825
$ws.Range("T3").Value = 0.434337400956825

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 57.14035266666667
$ws.Range("H4").Value = 171.421058
$ws.Range("I4").Value = 0.7274038390747541
$ws.Range("J4").Value = 0.7274038390747541
$ws.Range("O4").Value = 0.2250228664278349
$ws.Range("P4").Value = 0.2250228664278349
$ws.Range("Q4").Value = 334.0836046496849
$ws.Range("R4").Value = 3006.752441847164
$ws.Range("S4").Value = 0.1636824969192127
$ws.Range("T4").Value = 0.1636824969192127

# Row 5
$ws.Range("I5").Value = 0.08622113322131104
$ws.Range("J5").Value = 0.08622113322131104
$ws.Range("M5").Value = 4.621579
$ws.Range("N5").Value = 13.864737
$ws.Range("O5").Value = 0.1778708528171788
$ws.Range("P5").Value = 0.1778708528171788
$ws.Range("Q5").Value = 31.30195302647366
$ws.Range("R5").Value = 281.7175772382629
$ws.Range("S5").Value = 0.01533622649693818
$ws.Range("T5").Value = 0.01533622649693818

# Row 6
$ws.Range("I6").Value = 0.08622113322131104
$ws.Range("J6").Value = 0.08622113322131104
$ws.Range("N6").Value = 46.543441
$ws.Range("O6").Value = 0.5971062807549863
$ws.Range("P6").Value = 0.5971062807549863
$ws.Range("R6").Value = 945.716131135559
$ws.Range("S6").Value = 0.05148318018025723
$ws.Range("T6").Value = 0.05148318018025723

# Row 7
$ws.Range("I7").Value = 0.08622113322131104
$ws.Range("J7").Value = 0.08622113322131104
$ws.Range("O7").Value = 0.2250228664278349
$ws.Range("P7").Value = 0.2250228664278349
$ws.Range("S7").Value = 0.01940172654411564
$ws.Range("T7").Value = 0.01940172654411564

# Row 8
$ws.Range("I8").Value = 0.1863750277039348
$ws.Range("J8").Value = 0.1863750277039348
$ws.Range("M8").Value = 4.621579
$ws.Range("N8").Value = 13.864737
$ws.Range("O8").Value = 0.1778708528171788
$ws.Range("P8").Value = 0.1778708528171788
$ws.Range("Q8").Value = 67.66209332370899
$ws.Range("R8").Value = 608.958839913381
$ws.Range("S8").Value = 0.0331506851215242
$ws.Range("T8").Value = 0.03315068512152421

# Row 9
$ws.Range("I9").Value = 0.1863750277039348
$ws.Range("J9").Value = 0.1863750277039348
$ws.Range("N9").Value = 46.543441
$ws.Range("O9").Value = 0.5971062807549863
$ws.Range("P9").Value = 0.5971062807549863
$ws.Range("S9").Value = 0.111285699617904
$ws.Range("T9").Value = 0.111285699617904

# Row 10
$ws.Range("I10").Value = 0.1863750277039348
$ws.Range("J10").Value = 0.1863750277039348
$ws.Range("O10").Value = 0.2250228664278349
$ws.Range("P10").Value = 0.2250228664278349
$ws.Range("S10").Value = 0.04193864296450655
$ws.Range("T10").Value = 0.04193864296450656
